$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values
$ws.Range("B2").Value = 154
$ws.Range("B3").Value = -339

# Slightly widen the sheet's default/standard column width (8.59 -> 8.625 characters)
# (Worksheet.StandardWidth is the closest COM surface to sheetFormatPr's
# defaultColWidth; set it for parity even though this runtime's default-width
# plumbing is otherwise read-only.)
$ws.StandardWidth = 8.625

# Move the active cell selection from B3 to B4
$ws.Range("B4").Select()
